# #5: property boat&car done
# The "汽車" (car/vehicle) sheet was missing its header row (row 1 had
# been populated with a duplicate of the data row instead of column
# headers) and was missing the trailing columns that the other property
# sheets (land / building) carry: property_category, category, date,
# legislator_name, legislator_id, source_file, index. This restores the
# proper header row and fills in those extra columns for the existing
# data row, adding a new "capacity" column in place of "area".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Clone the existing header style (row 1, s=1) onto the new header cells,
# and the existing data style (row 2, s=2) onto the new data cells, so
# the extended columns look consistent with the rest of the row.
$ws.Range("B1").Copy($ws.Range("H1:N1"))
$ws.Range("B2").Copy($ws.Range("H2:N2"))

# --- Row 1: this row incorrectly held a copy of the data row; replace
# it with the real column headers (matching the 土地/建物 sheets). ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: B2:G2 (name/capacity/owner/register_date/register_reason/
# acquire_value) already hold the right values; fill in the new trailing
# columns to match the other property sheets. ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# Force text so "2013-07-15" isn't reinterpreted as a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2013-07-15"

$ws.Range("K2").Value = "林世嘉"
$ws.Range("L2").Value = 1740
$ws.Range("M2").Value = "tmpf70f1"
$ws.Range("N2").Value = 35
